# Apply updated profit calculation values to multiple Leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2365.375
$ws.Range("I31").Value = 224.6
$ws.Range("J31").Value = 5933.3335
$ws.Range("K31").Value = 673.8
$ws.Range("L31").Value = 17800.0005
$ws.Range("M31").Value = -443.8
$ws.Range("N31").Value = -18260.0005
$ws.Range("H70").Value = 1453.1818
$ws.Range("I70").Value = 802
$ws.Range("J70").Value = 1518.3
$ws.Range("K70").Value = 2406
$ws.Range("L70").Value = 4554.9
$ws.Range("M70").Value = -2136
$ws.Range("N70").Value = -5094.9
$ws.Range("H73").Value = 1453.1818
$ws.Range("I73").Value = 802
$ws.Range("J73").Value = 1518.3
$ws.Range("K73").Value = 2406
$ws.Range("L73").Value = 4554.9
$ws.Range("M73").Value = -1470
$ws.Range("N73").Value = -6426.9
$ws.Range("H82").Value = 877.4286
$ws.Range("I82").Value = 877.4286
$ws.Range("K82").Value = 2632.2858
$ws.Range("M82").Value = -2226.2858
$ws.Range("H85").Value = 877.4286
$ws.Range("I85").Value = 877.4286
$ws.Range("K85").Value = 2632.2858
$ws.Range("M85").Value = -1228.2858
$ws.Range("H129").Value = 1028.1428
$ws.Range("J129").Value = 1166.6666
$ws.Range("L129").Value = 3499.9998
$ws.Range("N129").Value = -13499.9998
$ws.Range("H141").Value = 1414
$ws.Range("I141").Value = 805.2632
$ws.Range("J141").Value = 2465.4546
$ws.Range("K141").Value = 2415.7896
$ws.Range("L141").Value = 7396.3638
$ws.Range("M141").Value = 2764.2104
$ws.Range("N141").Value = -17756.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4425.7417
$ws.Range("I32").Value = 2865.0532
$ws.Range("K32").Value = 2865.0532
$ws.Range("M32").Value = -2578.0532
$ws.Range("H132").Value = 1593.7869
$ws.Range("I132").Value = 1640.6863
$ws.Range("J132").Value = 1354.6
$ws.Range("K132").Value = 4922.0589
$ws.Range("L132").Value = 4063.8
$ws.Range("M132").Value = -2392.0589
$ws.Range("N132").Value = -9123.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1411.5
$ws.Range("I20").Value = 1186.6666
$ws.Range("J20").Value = 1718.091
$ws.Range("K20").Value = 1186.6666
$ws.Range("L20").Value = 1718.091
$ws.Range("M20").Value = -939.6666
$ws.Range("N20").Value = -2212.091
$ws.Range("H64").Value = 290
$ws.Range("I64").Value = 350
$ws.Range("J64").Value = 250
$ws.Range("K64").Value = 350
$ws.Range("L64").Value = 250
$ws.Range("M64").Value = -125
$ws.Range("N64").Value = -700
$ws.Range("H67").Value = 290
$ws.Range("I67").Value = 350
$ws.Range("J67").Value = 250
$ws.Range("K67").Value = 350
$ws.Range("L67").Value = 250
$ws.Range("M67").Value = 430
$ws.Range("N67").Value = -1810
$ws.Range("H86").Value = 1637.3611
$ws.Range("I86").Value = 1568.909
$ws.Range("J86").Value = 1744.9286
$ws.Range("K86").Value = 1568.909
$ws.Range("L86").Value = 1744.9286
$ws.Range("M86").Value = -445.9090000000001
$ws.Range("N86").Value = -3990.9286
$ws.Range("H89").Value = 1637.3611
$ws.Range("I89").Value = 1568.909
$ws.Range("J89").Value = 1744.9286
$ws.Range("K89").Value = 7844.545
$ws.Range("L89").Value = 8724.643
$ws.Range("M89").Value = -2228.545
$ws.Range("N89").Value = -19956.643
$ws.Range("H107").Value = 2727.75
$ws.Range("I107").Value = 2727.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2727.75
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -807.75
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 731351.5
$ws.Range("I134").Value = 1337123.6
$ws.Range("J134").Value = 4425
$ws.Range("K134").Value = 4011370.8
$ws.Range("L134").Value = 13275
$ws.Range("M134").Value = -4008835.8
$ws.Range("N134").Value = -18345

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 468.66666
$ws.Range("I47").Value = 103
$ws.Range("J47").Value = 1200
$ws.Range("K47").Value = 309
$ws.Range("L47").Value = 3600
$ws.Range("M47").Value = 122
$ws.Range("N47").Value = -4462
$ws.Range("H113").Value = 1212615.4
$ws.Range("I113").Value = 1894435.2
$ws.Range("J113").Value = 491.22223
$ws.Range("K113").Value = 5683305.6
$ws.Range("L113").Value = 1473.66669
$ws.Range("M113").Value = -5681135.6
$ws.Range("N113").Value = -5813.66669
$ws.Range("H115").Value = 2653.7144
$ws.Range("I115").Value = 1028
$ws.Range("J115").Value = 3304
$ws.Range("K115").Value = 3084
$ws.Range("L115").Value = 9912
$ws.Range("M115").Value = -1909
$ws.Range("N115").Value = -12262
$ws.Range("H131").Value = 920.91
$ws.Range("I131").Value = 900
$ws.Range("J131").Value = 921.1212
$ws.Range("K131").Value = 2700
$ws.Range("L131").Value = 2763.3636
$ws.Range("M131").Value = 2340
$ws.Range("N131").Value = -12843.3636
$ws.Range("H133").Value = 4631.905
$ws.Range("J133").Value = 6377.75
$ws.Range("L133").Value = 19133.25
$ws.Range("N133").Value = -29253.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 758.7917
$ws.Range("I16").Value = 713.5217
$ws.Range("J16").Value = 1800
$ws.Range("K16").Value = 713.5217
$ws.Range("L16").Value = 1800
$ws.Range("M16").Value = -543.5217
$ws.Range("N16").Value = -2140
$ws.Range("H61").Value = 2340
$ws.Range("I61").Value = 1900
$ws.Range("J61").Value = 2780
$ws.Range("K61").Value = 1900
$ws.Range("L61").Value = 2780
$ws.Range("M61").Value = -1698
$ws.Range("N61").Value = -3184
$ws.Range("H113").Value = 2340
$ws.Range("I113").Value = 1900
$ws.Range("J113").Value = 2780
$ws.Range("K113").Value = 1900
$ws.Range("L113").Value = 2780
$ws.Range("M113").Value = 270
$ws.Range("N113").Value = -7120
$ws.Range("H136").Value = 1564.4865
$ws.Range("I136").Value = 1217.1538
$ws.Range("J136").Value = 2385.4546
$ws.Range("K136").Value = 3651.4614
$ws.Range("L136").Value = 7156.3638
$ws.Range("M136").Value = -1101.4614
$ws.Range("N136").Value = -12256.3638
$ws.Range("H141").Value = 161750
$ws.Range("J141").Value = 161750
$ws.Range("L141").Value = 161750
$ws.Range("N141").Value = -172110

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2237.5
$ws.Range("I62").Value = 2800
$ws.Range("J62").Value = 1900
$ws.Range("K62").Value = 2800
$ws.Range("L62").Value = 1900
$ws.Range("M62").Value = -2176
$ws.Range("N62").Value = -3148
$ws.Range("H65").Value = 2237.5
$ws.Range("I65").Value = 2800
$ws.Range("J65").Value = 1900
$ws.Range("K65").Value = 14000
$ws.Range("L65").Value = 9500
$ws.Range("M65").Value = -10880
$ws.Range("N65").Value = -15740
